# Daily attendance processing - 2025-10-17 06:16:02
# Re-orders the "Recorded By" (column G) entries on the "Session Analysis Results" sheet
# by rotating the recorder list left by one (moving the first recorder to the end).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$recordedByUpdates = @{
    'G2' = 'system, backup@backdoor.com, System'
    'G3' = 'dnasr281@gmail.com, System'
    'G4' = 'backup@backdoor.com, System'
    'G5' = 'backup@backdoor.com, System'
    'G6' = 'dnasr281@gmail.com, System'
    'G10' = 'dnasr281@gmail.com, System'
    'G11' = 'dnasr281@gmail.com, System'
    'G12' = 'dnasr281@gmail.com, System'
    'G13' = 'dnasr281@gmail.com, System'
    'G14' = 'dnasr281@gmail.com, System'
    'G15' = 'dnasr281@gmail.com, System'
    'G17' = 'dnasr281@gmail.com, System'
    'G18' = 'dnasr281@gmail.com, System'
    'G19' = 'dnasr281@gmail.com, System'
    'G29' = 'system, backup@backdoor.com, System'
    'G30' = 'dnasr281@gmail.com, System'
    'G31' = 'backup@backdoor.com, System'
    'G32' = 'backup@backdoor.com, System'
    'G33' = 'dnasr281@gmail.com, System'
    'G37' = 'dnasr281@gmail.com, System'
    'G38' = 'dnasr281@gmail.com, System'
    'G39' = 'dnasr281@gmail.com, System'
    'G40' = 'dnasr281@gmail.com, System'
    'G41' = 'dnasr281@gmail.com, System'
    'G42' = 'dnasr281@gmail.com, System'
    'G44' = 'dnasr281@gmail.com, System'
    'G45' = 'dnasr281@gmail.com, System'
    'G46' = 'dnasr281@gmail.com, System'
    'G56' = 'system, backup@backdoor.com, System'
    'G57' = 'dnasr281@gmail.com, System'
    'G58' = 'backup@backdoor.com, System'
    'G59' = 'backup@backdoor.com, System'
    'G60' = 'dnasr281@gmail.com, System'
    'G64' = 'dnasr281@gmail.com, System'
    'G65' = 'dnasr281@gmail.com, System'
    'G66' = 'dnasr281@gmail.com, System'
    'G67' = 'dnasr281@gmail.com, System'
    'G68' = 'dnasr281@gmail.com, System'
    'G69' = 'dnasr281@gmail.com, System'
    'G71' = 'dnasr281@gmail.com, System'
    'G72' = 'dnasr281@gmail.com, System'
    'G73' = 'dnasr281@gmail.com, System'
    'G83' = 'backup@backdoor.com, System'
    'G84' = 'backup@backdoor.com, System'
    'G85' = 'backup@backdoor.com, System'
    'G86' = 'dnasr281@gmail.com, System'
    'G87' = 'dnasr281@gmail.com, System'
    'G88' = 'dnasr281@gmail.com, System'
    'G89' = 'dnasr281@gmail.com, System'
    'G90' = 'dnasr281@gmail.com, admin@admin.com'
    'G93' = 'dnasr281@gmail.com, System'
    'G95' = 'dnasr281@gmail.com, System'
    'G96' = 'dnasr281@gmail.com, System'
    'G97' = 'dnasr281@gmail.com, System'
    'G99' = 'dnasr281@gmail.com, System'
    'G109' = 'backup@backdoor.com, System'
    'G110' = 'backup@backdoor.com, System'
    'G111' = 'backup@backdoor.com, System'
    'G112' = 'dnasr281@gmail.com, System'
    'G113' = 'dnasr281@gmail.com, System'
    'G114' = 'dnasr281@gmail.com, System'
    'G115' = 'dnasr281@gmail.com, System'
    'G116' = 'dnasr281@gmail.com, admin@admin.com'
    'G119' = 'dnasr281@gmail.com, System'
    'G121' = 'dnasr281@gmail.com, System'
    'G122' = 'dnasr281@gmail.com, System'
    'G123' = 'dnasr281@gmail.com, System'
    'G125' = 'dnasr281@gmail.com, System'
    'G135' = 'backup@backdoor.com, System'
    'G136' = 'backup@backdoor.com, System'
    'G137' = 'backup@backdoor.com, System'
    'G138' = 'dnasr281@gmail.com, System'
    'G139' = 'dnasr281@gmail.com, System'
    'G140' = 'dnasr281@gmail.com, System'
    'G141' = 'dnasr281@gmail.com, System'
    'G142' = 'dnasr281@gmail.com, admin@admin.com'
    'G145' = 'dnasr281@gmail.com, System'
    'G147' = 'dnasr281@gmail.com, System'
    'G148' = 'dnasr281@gmail.com, System'
    'G149' = 'dnasr281@gmail.com, System'
    'G151' = 'dnasr281@gmail.com, System'
}

foreach ($cell in $recordedByUpdates.Keys) {
    $ws.Range($cell).Value = $recordedByUpdates[$cell]
}
